$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list: refresh Price (D) / Volume(1h) (E) columns for this pull.
# These are stored as literal text (the source feed writes them as strings, e.g.
# "246.16" / "0.27%"), so format each target cell as Text before assigning the
# value to stop Excel from auto-converting the numeric-looking text into a number
# or percentage value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.16"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.27%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.75%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.153"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.29%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05804"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.50%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.659"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.48%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.234"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "7.30%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8524"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.40%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8641"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.39%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1378"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.36%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07083"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.56%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03207"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "12.10%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09381"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.11%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001541"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.26%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006017"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-94.08%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005937"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.93%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.498"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.44%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.13%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3195"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03367"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.38%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.70%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.490"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.99%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04150"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.59%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1380"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.47%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001225"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.19%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.50%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001100"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.45%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "4.21%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03748"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.79%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005798"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.52%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1070"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.48%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.13%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008643"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-13.81%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005290"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.26%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.35%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-35.35%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002180"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-20.98%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.35%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.35%"

Write-Output "Updated cryptos price/volume values on $($ws.Name)"
